$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; existing rows 34:123 shift down to 35:124,
# carrying their data/formatting (including the date style on column D) with them.
$ws.Rows(34).Insert()

# Populate the newly inserted row 34 with the new record's data.
$ws.Cells.Item(34, 1).Value = 8
$ws.Cells.Item(34, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44614
$ws.Cells.Item(34, 5).Value = 4
$ws.Cells.Item(34, 6).Value = 100112040
$ws.Cells.Item(34, 7).Value = "Cilantro"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 1800
$ws.Cells.Item(34, 11).Value = 2300
$ws.Cells.Item(34, 12).Value = 2500
$ws.Cells.Item(34, 13).Value = 2400
$ws.Cells.Item(34, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(34, 16).Value = 1600
$ws.Cells.Item(34, 17).Value = 1.5
$ws.Cells.Item(34, 18).Value = "Hortaliza"
